# Update "想去人数" (F column) figures for the first three events on the
# "展览" (Exhibition) and "全部类型" (All types) sheets, matching the
# newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 9243
    $ws.Range("F3").Value = 207
    $ws.Range("F4").Value = 489
}
